# Re-process the data with the newly curated dimensions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: metadata "kind" strings - fix mislabeled dimension -> measure,
# and switch the "aragon" column to the sdmx reference-area dimension.
$ws.Range("A2").Value = "iaest-measure:horas-trabajadas"
$ws.Range("D2").Value = "sdmx-dimension:refArea"

# Row 3: "horas-trabajadas" was wrongly curated as a dimension ("dim"); it is
# actually a measure ("medida"), matching columns B and C.
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "medida"
$ws.Range("D3").Value = "dim"

# Row 4: data types follow the corrected roles above.
$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("D4").Value = "URI-Comunidad"

# Row 5 (the old per-column mapping-file references) is no longer needed now
# that the dimensions are curated directly, so drop the whole row.
$ws.Range("A5:E5").Delete()
